# Update the cryptocurrency price/volume table with freshly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates -------------------------------------------
# Values that parse as plain numbers need a leading apostrophe so Excel
# keeps them stored as text (matching the sheet's existing text layout)
# instead of silently converting them to numeric cells.
$ws.Range("D2").Value  = "27.881.14"
$ws.Range("D3").Value  = "1.634.96"
$ws.Range("D5").Value  = "'211.76"
$ws.Range("D8").Value  = "'23.39"
$ws.Range("D11").Value = "'0.0884"
$ws.Range("D12").Value = "1.866.14"
$ws.Range("D13").Value = "1.640.76"
$ws.Range("D16").Value = "'65.29"
$ws.Range("D17").Value = "27.884.07"
$ws.Range("D18").Value = "'229.06"
$ws.Range("D22").Value = "'4.34"
$ws.Range("D23").Value = "'10.03"
$ws.Range("D24").Value = "'2.07"
$ws.Range("D25").Value = "'155.26"
$ws.Range("D28").Value = "'15.53"
$ws.Range("D34").Value = "1.394.94"
$ws.Range("D37").Value = "'2.34"
$ws.Range("D45").Value = "'5.44"
$ws.Range("D46").Value = "1.774.38"
$ws.Range("D48").Value = "'88.73"

# --- Column E (Volume 1h) updates ---------------------------------------
$ws.Range("E2").Value  = "  -0.08%  "
$ws.Range("E3").Value  = "  -0.03%  "
$ws.Range("E4").Value  = "  +0.03%  "
$ws.Range("E5").Value  = "  -0.23%  "
$ws.Range("E6").Value  = "  -0.54%  "
$ws.Range("E7").Value  = "  +0.03%  "
$ws.Range("E8").Value  = "  +0.61%  "
$ws.Range("E9").Value  = "  -0.58%  "
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("E19").Value = "  +2.44%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("E23").Value = "  -2.81%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("E35").Value = "  +1.56%  "
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("E40").Value = "  -2.84%  "
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("E45").Value = "  -1.50%  "
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("E47").Value = "  -2.59%  "
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("E51").Value = "  +0.50%  "

# --- Rows 43 and 44: RenderToken and Aave swapped places in the ranking --
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'1.83"
$ws.Range("E43").Value = "  +1.51%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'65.89"
$ws.Range("E44").Value = "  -1.54%  "
